# Edit script: update struggle dataset values (C2:H21), remove now-obsolete row 22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for columns C:H across rows 2-21 (ax, ay, az, gx, gy, gz)
$newData = @(
    @(3.09836453199387,-3.448286980390552,1.224327325820929,0.4437935948371887,-0.864068865776062,-0.9851729273796082),
    @(5.926505327224739,-2.135328322648994,-0.6819987297058234,1.0256427526474,-4.58713960647583,0.0221438650041818),
    @(-1.439981520175948,-1.973837949335578,-4.723978996276856,2.594497442245483,0.8452847599983215,-2.786614418029785),
    @(-13.93973755836486,-5.416623592376709,-4.78411507606506,2.057699680328369,1.456302762031555,-0.2620611786842346),
    @(-5.790439605712891,-5.474406182765961,-2.479388117790221,-0.5109887719154358,-1.112385630607605,-1.43751859664917),
    @(-6.174473762512195,-4.646332740783681,-6.182556152343746,-0.303600013256073,2.510961532592773,-0.3182607889175415),
    @(-1.106817305088029,0.3479279279709004,-4.026142060756678,-0.5890268087387085,-0.6258314251899719,0.7171558141708374),
    @(-0.1755727529525754,2.31334447860718,-5.915446519851686,1.278999090194702,0.4355469346046448,0.9905179738998412),
    @(1.180307447910312,0.6379154920577954,-5.021014630794521,-1.44057297706604,-0.845132052898407,-0.7629706859588623),
    @(-1.093243360519415,-1.983224630355838,-1.879333615303035,-1.023504734039307,-0.6563746929168701,0.645684540271759),
    @(-3.451674103736873,-1.611639708280552,-0.6418006122112334,-1.10917854309082,-0.5216789245605469,0.8017606139183044),
    @(-2.273676156997681,0.1646193265914786,-3.850346922874469,-3.361896991729736,3.937331914901733,2.076178312301636),
    @(-2.109282225370405,-3.67548027634622,-7.529284000396693,2.992323398590088,0.9390525817871094,0.3888157308101654),
    @(-1.644850492477401,-7.059904575347879,1.617063522338856,-0.3591887652873993,1.687973380088806,0.8848382830619812),
    @(1.835043907165529,-2.18121553957462,-0.952013134956361,1.22615921497345,0.5057964324951172,0.3527746796607971),
    @(8.433930218219764,-3.52629014849663,2.835070371627813,-0.1565342247486114,-0.5474879741668701,0.3535382747650146),
    @(-5.246673464775119,-3.962655484676361,-5.154342770576502,0.4882340431213379,-7.02800464630127,-4.413654327392578),
    @(-7.417413711547821,-2.038821458816521,-6.523755788803086,1.305877208709717,-2.09221339225769,-1.915520668029785),
    @(0.6842400431633051,-3.085346877574939,-2.034696012735365,1.73409366607666,-3.325092315673828,-2.228436470031738),
    @(1.740720510482789,-5.68582010269165,-1.830426752567291,-2.827542543411255,-0.811687171459198,-1.57450520992279)
)

$startRow = 2
for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $startRow + $i
    $vals = $newData[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = 3 + $c  # column C = 3
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}

# Remove row 22, which is no longer part of the dataset
$ws.Rows.Item(22).Delete()
